# Updated cryptos list — applies the latest scraped price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price" (column D) cell while forcing it to
# remain plain text, since many prices look like numbers (or even
# multi-dot "numbers" such as 34.189.87) and must not be silently coerced
# into a numeric type by Excel's auto-detection. Temporarily flipping the
# cell to Text format for the assignment and then clearing formats again
# keeps the cell's original (default) style intact.
function Set-PriceText {
    param($addr, $text)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# --- Row 2: Bitcoin ---
Set-PriceText "D2" "34.189.87"
$ws.Range("E2").Value = "  +0.33%  "

# --- Row 3: Ethereum ---
Set-PriceText "D3" "1.788.99"
$ws.Range("E3").Value = "  -0.02%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.18%  "

# --- Row 5: BNB ---
Set-PriceText "D5" "226.17"
$ws.Range("E5").Value = "  -0.38%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +0.58%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.15%  "

# --- Row 8: Solana ---
Set-PriceText "D8" "32.37"
$ws.Range("E8").Value = "  +0.39%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +0.16%  "

# --- Row 10: Dogecoin ---
Set-PriceText "D10" "0.0689"
$ws.Range("E10").Value = "  +0.28%  "

# --- Row 11: TRON ---
Set-PriceText "D11" "0.0947"
$ws.Range("E11").Value = "  +0.67%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
Set-PriceText "D12" "2.047.56"
$ws.Range("E12").Value = "  +0.01%  "

# --- Row 13: Chainlink ---
Set-PriceText "D13" "11.20"
$ws.Range("E13").Value = "  -1.31%  "

# --- Row 14: WrappedEther ---
Set-PriceText "D14" "1.786.81"
$ws.Range("E14").Value = "  +0.10%  "

# --- Row 15: Polygon ---
Set-PriceText "D15" "0.626"
$ws.Range("E15").Value = "  +0.54%  "

# --- Row 16: WrappedBTC ---
Set-PriceText "D16" "34.183.80"
$ws.Range("E16").Value = "  +0.32%  "

# --- Row 17: Polkadot ---
$ws.Range("E17").Value = "  +0.46%  "

# --- Row 19: ShibaInu ---
Set-PriceText "D19" "0.0₃0806"
$ws.Range("E19").Value = "  +3.16%  "

# --- Row 20: BitcoinCash ---
Set-PriceText "D20" "246.17"
$ws.Range("E20").Value = "  +0.98%  "

# --- Row 21: Avalanche ---
Set-PriceText "D21" "11.04"
$ws.Range("E21").Value = "  +0.89%  "

# --- Row 22: Dai ---
$ws.Range("E22").Value = "  +0.20%  "

# --- Row 23: Uniswap ---
Set-PriceText "D23" "4.17"
$ws.Range("E23").Value = "  +1.79%  "

# --- Row 24: Toncoin ---
$ws.Range("E24").Value = "  +0.58%  "

# --- Row 25: Monero ---
Set-PriceText "D25" "161.85"
$ws.Range("E25").Value = "  -0.10%  "

# --- Row 26: Cosmos ---
Set-PriceText "D26" "7.17"
$ws.Range("E26").Value = "  -0.38%  "

# --- Row 27: EthereumClassic ---
Set-PriceText "D27" "16.31"
$ws.Range("E27").Value = "  +0.27%  "

# --- Row 28: Stellar ---
$ws.Range("E28").Value = "  +1.04%  "

# --- Row 29: BinanceUSD ---
Set-PriceText "D29" "1.00"
$ws.Range("E29").Value = "  +0.24%  "

# --- Row 30: Hedera ---
$ws.Range("E30").Value = "  +0.08%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  -0.64%  "

# --- Row 32: Filecoin ---
Set-PriceText "D32" "3.77"
$ws.Range("E32").Value = "  +3.37%  "

# --- Row 34: LidoDAOToken ---
$ws.Range("E34").Value = "  -1.47%  "

# --- Row 35: Maker ---
Set-PriceText "D35" "1.444.53"
$ws.Range("E35").Value = "  +2.36%  "

# --- Row 36: RenderToken ---
Set-PriceText "D36" "2.59"
$ws.Range("E36").Value = "  +9.40%  "

# --- Row 37: ImmutableX ---
Set-PriceText "D37" "0.665"
$ws.Range("E37").Value = "  +2.82%  "

# --- Rows 38 & 39: VeChain and TrustWalletToken swapped rank order ---
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-PriceText "D38" "1.05"
$ws.Range("E38").Value = "  +1.12%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-PriceText "D39" "0.0191"
$ws.Range("E39").Value = "  +0.04%  "

# --- Row 40: Aave ---
Set-PriceText "D40" "82.17"
$ws.Range("E40").Value = "  +2.11%  "

# --- Row 42: ARBITRUM ---
Set-PriceText "D42" "0.923"
$ws.Range("E42").Value = "  +0.22%  "

# --- Row 43: MXToken ---
$ws.Range("E43").Value = "  +1.08%  "

# --- Row 44: InjectiveProtocol ---
Set-PriceText "D44" "13.71"
$ws.Range("E44").Value = "  +2.79%  "

# --- Row 45: Kaspa ---
Set-PriceText "D45" "0.0519"
$ws.Range("E45").Value = "  +2.30%  "

# --- Row 46: FraxShare ---
Set-PriceText "D46" "6.12"
$ws.Range("E46").Value = "  +1.40%  "

# --- Row 47: WEMIXToken ---
$ws.Range("E47").Value = "  +0.74%  "

# --- Row 48: RocketPoolETH ---
Set-PriceText "D48" "1.946.38"
$ws.Range("E48").Value = "  -0.07%  "

# --- Row 49: Quant ---
Set-PriceText "D49" "105.24"
$ws.Range("E49").Value = "  -1.50%  "

# --- Row 50: PaxDollar ---
$ws.Range("E50").Value = "  +0.17%  "

# --- Row 51: BabyDogeCoin ---
Set-PriceText "D51" "0.0₆0129"
$ws.Range("E51").Value = "  -7.05%  "
